# Menu system + splash RLE: rework the PSC/ARR prescaler values, add the
# MIC column header next to the buffer layout, and add a "CW FFT width"
# row (mirroring the existing "SSB FFT width" row) above the PWM-bits row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- ADC channels (B2) and ARR divider (B4) changed ---
$ws.Range("B2").Value = 2
$ws.Range("B4").Value = 2

# --- MIC label moves from E3 to the new F1 header cell ---
$ws.Range("E3").ClearContents()
$ws.Range("F1").Value = "MIC"

# --- Insert a new row above the final (PWM bits) row, shifting it down ---
$ws.Rows(22).Insert()

# --- New row 21: CW FFT width = 150 / B18, labeled in "bins" ---
$ws.Range("A21").Value = "CW FFT width"
$ws.Range("B21").Formula = "=150/B18"
$ws.Range("C21").Value = "bins"

# Match the formatting used by the other computed cells (fill style s="1")
$ws.Range("B18").Copy()
$ws.Range("B21").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- Update the view: scrolled down a bit, new active selection ---
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("H15").Select()
